$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")
$ws.Activate()

# Copy the existing date cell formatting (C76:D76) onto the new date cells first,
# so the new rows reuse the workbook's existing date-format style (numFmtId 14)
# instead of Excel creating a brand-new number format.
$ws.Range("C76:D76").Copy()
$ws.Range("C77:D78").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 77: Talk Like Ted
$ws.Range("A77").Value = "Talk Like Ted"
$ws.Range("B77").Value = "Carmine Gallo"
$ws.Range("C77").Value = 43962
$ws.Range("D77").Value = 43965
$ws.Range("E77").Value = "presentation;speaking;persuasion;presence"
$ws.Range("F77").Value = "Audio"
$ws.Range("G77").Value = "7 Hours 42 Mins"

# Row 78: The Elephant in the Brain
$ws.Range("A78").Value = "The Elephant in the Brain"
$ws.Range("B78").Value = "Kevin Simler"
$ws.Range("C78").Value = 43965
$ws.Range("D78").Value = 43970
$ws.Range("E78").Value = "psycology;sociology;brain;humans;politics"
$ws.Range("F78").Value = "Audio"
$ws.Range("G78").Value = "10 Hours 32 Mins"

# Update the view so the new rows are visible, matching the saved workbook state
$excel.ActiveWindow.ScrollRow = 59
$ws.Range("A79").Select()
